$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 2.8
$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 3.5
$ws.Range("I2").Value = 3.7
$ws.Range("J2").Value = 2.7
$ws.Range("K2").Value = 2.72
$ws.Range("L2").Value = 2.76
$ws.Range("M2").Value = 1.27
$ws.Range("N2").Value = 1.8
$ws.Range("O2").Value = 2.18
$ws.Range("P2").Value = 1.23
$ws.Range("Q2").Value = 5.1
$ws.Range("R2").Value = 1.06
$ws.Range("S2").Value = 15
$ws.Range("T2").Value = 3.3
$ws.Range("U2").Value = 1.39
$ws.Range("V2").Value = 1.37
$ws.Range("W2").Value = 1.51
$ws.Range("X2").Value = 4.6
$ws.Range("Y2").Value = 7
$ws.Range("Z2").Value = 23
$ws.Range("AA2").Value = 130
$ws.Range("AB2").Value = 5.7
$ws.Range("AD2").Value = 23
$ws.Range("AE2").Value = 140
$ws.Range("AF2").Value = 15.5
$ws.Range("AG2").Value = 23
$ws.Range("AH2").Value = 55
$ws.Range("AI2").Value = 260
$ws.Range("AJ2").Value = 90
$ws.Range("AK2").Value = 110
$ws.Range("AL2").Value = 250
$ws.Range("AN2").Value = 180
$ws.Range("AO2").Value = 270
$ws.Range("F3").Value = 2.4
$ws.Range("G3").Value = 2.46
$ws.Range("H3").Value = 4.7
$ws.Range("I3").Value = 4.9
$ws.Range("J3").Value = 2.6
$ws.Range("K3").Value = 2.66
$ws.Range("L3").Value = 2.72
$ws.Range("M3").Value = 1.27
$ws.Range("N3").Value = 1.83
$ws.Range("O3").Value = 2.16
$ws.Range("P3").Value = 1.23
$ws.Range("Q3").Value = 5.1
$ws.Range("R3").Value = 1.06
$ws.Range("S3").Value = 14.5
$ws.Range("T3").Value = 3.4
$ws.Range("U3").Value = 1.39
$ws.Range("V3").Value = 1.26
$ws.Range("X3").Value = 4.8
$ws.Range("Y3").Value = 8.4
$ws.Range("Z3").Value = 34
$ws.Range("AA3").Value = 200
$ws.Range("AB3").Value = 5
$ws.Range("AD3").Value = 29
$ws.Range("AE3").Value = 170
$ws.Range("AF3").Value = 11.5
$ws.Range("AG3").Value = 17
$ws.Range("AH3").Value = 65
$ws.Range("AI3").Value = 380
$ws.Range("AJ3").Value = 48
$ws.Range("AK3").Value = 75
$ws.Range("AL3").Value = 240
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 130
$ws.Range("AO3").Value = 530
$ws.Range("H4").Value = 24
$ws.Range("I4").Value = 27
$ws.Range("J4").Value = 10.5
$ws.Range("L4").Value = 1.18
$ws.Range("O4").Value = 1.08
$ws.Range("P4").Value = 4.2
$ws.Range("Q4").Value = 1.3
$ws.Range("R4").Value = 2.3
$ws.Range("S4").Value = 1.75
$ws.Range("T4").Value = 1.9
$ws.Range("U4").Value = 2.02
$ws.Range("W4").Value = 7.4
$ws.Range("X4").Value = 90
$ws.Range("Z4").Value = 1000
$ws.Range("AB4").Value = 19
$ws.Range("AC4").Value = 29
$ws.Range("AD4").Value = 980
$ws.Range("AE4").Value = 1000
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 14.5
$ws.Range("AH4").Value = 1000
$ws.Range("AI4").Value = 1000
$ws.Range("AJ4").Value = 11.5
$ws.Range("AK4").Value = 14.5
$ws.Range("AL4").Value = 42
$ws.Range("AM4").Value = 490
$ws.Range("AN4").Value = 2.56
$ws.Range("G5").Value = 2.04
$ws.Range("H5").Value = 5.3
$ws.Range("J5").Value = 3.05
$ws.Range("K5").Value = 3.2
$ws.Range("O5").Value = 1.68
$ws.Range("P5").Value = 1.39
$ws.Range("Q5").Value = 3.15
$ws.Range("W5").Value = 1.96
$ws.Range("AE5").Value = 180
$ws.Range("F6").Value = 2.64
$ws.Range("G6").Value = 2.66
$ws.Range("L6").Value = 1.63
$ws.Range("N6").Value = 2.72
$ws.Range("O6").Value = 1.57
$ws.Range("P6").Value = 1.56
$ws.Range("Q6").Value = 2.74
$ws.Range("R6").Value = 1.2
$ws.Range("S6").Value = 5.9
$ws.Range("T6").Value = 2.18
$ws.Range("U6").Value = 1.82
$ws.Range("X6").Value = 8.4
$ws.Range("Y6").Value = 8.800000000000001
$ws.Range("AB6").Value = 7.8
$ws.Range("AC6").Value = 6.6
$ws.Range("AI6").Value = 80
$ws.Range("AJ6").Value = 40
$ws.Range("AK6").Value = 38
$ws.Range("AM6").Value = 180
$ws.Range("AN6").Value = 42
$ws.Range("AO6").Value = 75
$ws.Range("F7").Value = 1.63
$ws.Range("G7").Value = 1.65
$ws.Range("I7").Value = 5.5
$ws.Range("J7").Value = 4.8
$ws.Range("K7").Value = 5.1
$ws.Range("L7").Value = 1.25
$ws.Range("N7").Value = 7.4
$ws.Range("O7").Value = 1.14
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.44
$ws.Range("R7").Value = 1.86
$ws.Range("S7").Value = 2.12
$ws.Range("T7").Value = 1.55
$ws.Range("U7").Value = 2.72
$ws.Range("V7").Value = 1.22
$ws.Range("W7").Value = 2.52
$ws.Range("X7").Value = 36
$ws.Range("Y7").Value = 36
$ws.Range("Z7").Value = 65
$ws.Range("AA7").Value = 1000
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 12.5
$ws.Range("AE7").Value = 55
$ws.Range("AF7").Value = 14
$ws.Range("AH7").Value = 17
$ws.Range("AI7").Value = 110
$ws.Range("AJ7").Value = 18
$ws.Range("AL7").Value = 23
$ws.Range("AM7").Value = 60
$ws.Range("F8").Value = 5.3
$ws.Range("G8").Value = 5.5
$ws.Range("H8").Value = 1.8
$ws.Range("I8").Value = 1.83
$ws.Range("N8").Value = 3.85
$ws.Range("O8").Value = 1.32
$ws.Range("T8").Value = 1.89
$ws.Range("V8").Value = 2.2
$ws.Range("W8").Value = 1.22
$ws.Range("Z8").Value = 10.5
$ws.Range("AA8").Value = 19.5
$ws.Range("AB8").Value = 18
$ws.Range("AC8").Value = 9
$ws.Range("AD8").Value = 10
$ws.Range("AE8").Value = 19
$ws.Range("AF8").Value = 42
$ws.Range("AG8").Value = 21
$ws.Range("AH8").Value = 21
$ws.Range("AJ8").Value = 140
$ws.Range("AK8").Value = 75
$ws.Range("AL8").Value = 80
$ws.Range("AM8").Value = 130
$ws.Range("AN8").Value = 85
$ws.Range("AO8").Value = 12.5
$ws.Range("F9").Value = 1.85
$ws.Range("G9").Value = 1.95
$ws.Range("H9").Value = 5.1
$ws.Range("I9").Value = 6.2
$ws.Range("J9").Value = 3.15
$ws.Range("K9").Value = 3.5
$ws.Range("L9").Value = 1.59
$ws.Range("N9").Value = 2.62
$ws.Range("O9").Value = 1.54
$ws.Range("P9").Value = 1.52
$ws.Range("Q9").Value = 2.66
$ws.Range("R9").Value = 1.18
$ws.Range("S9").Value = 5.5
$ws.Range("T9").Value = 2.26
$ws.Range("U9").Value = 1.66
$ws.Range("V9").Value = 1.2
$ws.Range("W9").Value = 2.04
$ws.Range("X9").Value = 9.199999999999999
$ws.Range("Y9").Value = 15.5
$ws.Range("Z9").Value = 50
$ws.Range("AB9").Value = 6.2
$ws.Range("AC9").Value = 8.4
$ws.Range("AD9").Value = 28
$ws.Range("AF9").Value = 10.5
$ws.Range("AG9").Value = 11.5
$ws.Range("AH9").Value = 30
$ws.Range("AJ9").Value = 25
$ws.Range("AK9").Value = 29
$ws.Range("AL9").Value = 70
$ws.Range("AN9").Value = 26
$ws.Range("F10").Value = 1.48
$ws.Range("G10").Value = 1.49
$ws.Range("H10").Value = 9.6
$ws.Range("I10").Value = 9.800000000000001
$ws.Range("J10").Value = 4.5
$ws.Range("K10").Value = 4.6
$ws.Range("L10").Value = 1.49
$ws.Range("M10").Value = 1.09
$ws.Range("N10").Value = 3.25
$ws.Range("O10").Value = 1.43
$ws.Range("P10").Value = 1.77
$ws.Range("Q10").Value = 2.26
$ws.Range("R10").Value = 1.27
$ws.Range("S10").Value = 4.3
$ws.Range("T10").Value = 2.46
$ws.Range("U10").Value = 1.64
$ws.Range("W10").Value = 3.05
$ws.Range("X10").Value = 12.5
$ws.Range("Y10").Value = 23
$ws.Range("Z10").Value = 80
$ws.Range("AA10").Value = 440
$ws.Range("AB10").Value = 5.8
$ws.Range("AD10").Value = 36
$ws.Range("AF10").Value = 7.2
$ws.Range("AK10").Value = 19
$ws.Range("AN10").Value = 10.5
$ws.Range("AO10").Value = 360
$ws.Range("F11").Value = 1.62
$ws.Range("G11").Value = 1.63
$ws.Range("P11").Value = 2.12
$ws.Range("Q11").Value = 1.87
$ws.Range("R11").Value = 1.43
$ws.Range("S11").Value = 3.2
$ws.Range("T11").Value = 1.92
$ws.Range("U11").Value = 2.06
$ws.Range("W11").Value = 2.58
$ws.Range("X11").Value = 16.5
$ws.Range("Y11").Value = 20
$ws.Range("Z11").Value = 48
$ws.Range("AA11").Value = 170
$ws.Range("AB11").Value = 8.4
$ws.Range("AC11").Value = 9.6
$ws.Range("AD11").Value = 23
$ws.Range("AE11").Value = 90
$ws.Range("AF11").Value = 9.199999999999999
$ws.Range("AH11").Value = 21
$ws.Range("AK11").Value = 15.5
$ws.Range("AM11").Value = 100
$ws.Range("AN11").Value = 8.800000000000001
$ws.Range("AO11").Value = 110
$ws.Range("F12").Value = 2.18
$ws.Range("G12").Value = 2.22
$ws.Range("I12").Value = 4.4
$ws.Range("L12").Value = 1.51
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 3.15
$ws.Range("O12").Value = 1.42
$ws.Range("Q12").Value = 2.3
$ws.Range("S12").Value = 4.3
$ws.Range("T12").Value = 1.89
$ws.Range("U12").Value = 1.98
$ws.Range("V12").Value = 1.3
$ws.Range("W12").Value = 1.81
$ws.Range("X12").Value = 10.5
$ws.Range("Z12").Value = 32
$ws.Range("AA12").Value = 110
$ws.Range("AD12").Value = 17.5
$ws.Range("AE12").Value = 65
$ws.Range("AH12").Value = 21
$ws.Range("AI12").Value = 75
$ws.Range("AJ12").Value = 27
$ws.Range("AM12").Value = 580
$ws.Range("AN12").Value = 23
$ws.Range("AO12").Value = 85
